$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the contents of F6 (it held the shared string "document"),
# while keeping the cell's existing formatting/style.
$ws.Range("F6").ClearContents()

# Move the active selection to F6 (was D11).
$ws.Activate()
$ws.Range("F6").Select()
